$d = $word.ActiveDocument

# Word Online couldn't cope with the numPr/jc and ilvl/numId ordering Word
# itself produced for list paragraphs inside table cells. Rewrite each of
# those paragraphs so that <w:numPr> (with <w:ilvl> before <w:numId>) comes
# before <w:jc> in <w:pPr>, matching the fixed-up canonical order.

function Set-ListParaOrder($para, $text, $numId) {
    $escaped = $text -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:pPr>' +
           '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="' + $numId + '"/></w:numPr>' +
           '<w:jc w:val="left"/>' +
           '</w:pPr><w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

Set-ListParaOrder $d.Paragraphs.Item(4) "Cell with"      "1001"
Set-ListParaOrder $d.Paragraphs.Item(5) "A"               "1001"
Set-ListParaOrder $d.Paragraphs.Item(6) "Bullet list"      "1001"
Set-ListParaOrder $d.Paragraphs.Item(7) "Cell with"        "1002"
Set-ListParaOrder $d.Paragraphs.Item(8) "A"                "1002"
Set-ListParaOrder $d.Paragraphs.Item(9) "Numbered list."    "1002"
